# Updated results of SP Class B
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row 6 (Export_Results moves from row 5 to row 6).
# Copy formatting from the row above (row 5) so the new A6 cell keeps
# the same style (border/bold/centered) as the other data cells in column A.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)  # xlPasteFormats

# Row 2: Call_Graph - time updated, name unchanged
$ws.Range("A2").Value = 89.65000000000001

# Row 3: new "Profiling" entry inserted (pushes following rows down)
$ws.Range("A3").Value = 92.53
$ws.Range("B3").Value = "Profiling"

# Row 4: Signal_Reconstruction (was row 3)
$ws.Range("A4").Value = 10.85
$ws.Range("B4").Value = "Signal_Reconstruction"

# Row 5: Energy_Estimation (was row 4)
$ws.Range("A5").Value = 2.41
$ws.Range("B5").Value = "Energy_Estimation"

# Row 6: Export_Results (was row 5)
$ws.Range("A6").Value = 19.38
$ws.Range("B6").Value = "Export_Results"

Write-Output "Updated results of SP Class B applied"
